$d = $word.ActiveDocument
$d.Content.Find.Execute("Expressing as rsin(θ±α) or rcos(θ±α)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Expressing as rsin(θ±α) or rcos(θ±α)", 2)
